# Handles float input without breaking stuff
#
# The "quiz" marksheet previously stored the student's per-question answers
# in three parallel blocks (A/B, D/E, G/H) because of a bug where results
# were duplicated across columns instead of being written once. This script
# fixes the summary numbers and collapses the answer-key table back down to
# a single clean block (with a small legacy remnant in D/E for the first
# three questions), matching the corrected export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Copy-CellStyle($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}

# ---------------------------------------------------------------------
# 1) Summary block (rows 10-12): fix values, give the row-label cells in
#    column A the same "mtitleStyle" formatting already used by A9.
# ---------------------------------------------------------------------
Copy-CellStyle "A9" "A10"
Copy-CellStyle "A9" "A11"
Copy-CellStyle "A9" "A12"

$ws.Range("B10").Value = 15
$ws.Range("C10").Value = 6
$ws.Range("D10").Value = 7
$ws.Range("E10").Value = 28

$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1   # was stored as the text "-1"; now a real number

$ws.Range("B12").Value = 60
$ws.Range("C12").Value = -6
$ws.Range("E12").Value = "54/112"   # was "Absent"

# ---------------------------------------------------------------------
# 2) Drop the G:H "third" answer block entirely, and the D:E block for
#    every question after the third (rows 19-40). A full .Clear() wipes
#    both content and formatting so the sheet's used range shrinks back
#    down from H40 to E40.
# ---------------------------------------------------------------------
$ws.Range("G15:H40").Clear() | Out-Null
$ws.Range("D19:E40").Clear() | Out-Null

# ---------------------------------------------------------------------
# 3) Populate the surviving D/E remnant (questions 1-3, rows 16-18) with
#    the student's (correct) answers for that duplicated block.
# ---------------------------------------------------------------------
Copy-CellStyle "B10" "D16"
$ws.Range("D16").Value = "Option A"

Copy-CellStyle "B10" "D17"
$ws.Range("D17").Value = "Option C"

Copy-CellStyle "B10" "D18"
$ws.Range("D18").Value = "Option D"

# ---------------------------------------------------------------------
# 4) Fill in the student's actual answers for the main A column
#    (questions 1-25, rows 16-40). Style mirrors correctness:
#      correctStyle (green)   -> matches the Correct Ans in column B
#      incorrectStyle (red)   -> differs from the Correct Ans in column B
#      normalStyle (black)    -> left blank = not attempted (no change needed)
# ---------------------------------------------------------------------
$studentAnswers = @{
    16 = "Option A"  # correct
    18 = "Option A"  # incorrect (correct ans is Option B)
    19 = "Option C"  # correct
    22 = "Option D"  # correct
    24 = "Option B"  # incorrect (correct ans is Option A)
    25 = "Option A"  # correct
    26 = "Option C"  # correct
    27 = "Option A"  # correct
    28 = "Option B"  # incorrect (correct ans is Option D)
    30 = "Option B"  # correct
    31 = "Option D"  # correct
    32 = "Option C"  # correct
    33 = "Option B"  # incorrect (correct ans is Option D)
    34 = "Option B"  # correct
    35 = "Option B"  # incorrect (correct ans is Option D)
    36 = "Option B"  # incorrect (correct ans is Option A)
    38 = "Option A"  # correct
    39 = "Option D"  # correct
}

$incorrectRows = @(18, 24, 28, 33, 35, 36)

foreach ($row in $studentAnswers.Keys) {
    $addr = "A$row"
    if ($incorrectRows -contains $row) {
        Copy-CellStyle "C10" $addr
    } else {
        Copy-CellStyle "B10" $addr
    }
    $ws.Range($addr).Value = $studentAnswers[$row]
}
